# Update the "想去人数" (want-to-go count) column F values on several rows
# across three worksheets: 展览 (Exhibitions), 演出 (Performances) and
# 全部类型 (All types), matching the upstream data refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 846
$ws1.Range("F3").Value  = 981
$ws1.Range("F4").Value  = 772
$ws1.Range("F5").Value  = 863
$ws1.Range("F6").Value  = 434
$ws1.Range("F7").Value  = 670
$ws1.Range("F9").Value  = 1266
$ws1.Range("F10").Value = 692
$ws1.Range("F11").Value = 407
$ws1.Range("F12").Value = 539
$ws1.Range("F13").Value = 179
$ws1.Range("F15").Value = 869
$ws1.Range("F16").Value = 7
$ws1.Range("F20").Value = 574
$ws1.Range("F21").Value = 135
$ws1.Range("F22").Value = 626
$ws1.Range("F24").Value = 931

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value  = 640
$ws2.Range("F7").Value  = 236
$ws2.Range("F10").Value = 26
$ws2.Range("F11").Value = 108

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 846
$ws4.Range("F5").Value  = 981
$ws4.Range("F6").Value  = 772
$ws4.Range("F7").Value  = 863
$ws4.Range("F8").Value  = 434
$ws4.Range("F9").Value  = 670
$ws4.Range("F11").Value = 1266
$ws4.Range("F12").Value = 692
$ws4.Range("F15").Value = 407
$ws4.Range("F16").Value = 539
$ws4.Range("F17").Value = 640
$ws4.Range("F18").Value = 179
$ws4.Range("F20").Value = 869
$ws4.Range("F22").Value = 7
$ws4.Range("F26").Value = 236
$ws4.Range("F28").Value = 574
$ws4.Range("F30").Value = 26
$ws4.Range("F31").Value = 108
$ws4.Range("F32").Value = 108
$ws4.Range("F33").Value = 135
$ws4.Range("F34").Value = 626
$ws4.Range("F36").Value = 931
